$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C81").Value = 0.55000000000000004
$ws.Range("C82").Value = 0.48
$ws.Range("C83").Value = 0.94
$ws.Range("C84").Value = 0.8
$ws.Range("C85").Value = 0.8
